$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2"=1.141508822445189; "D2"=0.1861054361842562; "E2"=0.1740423477911719; "F2"=1.317766729300061; "G2"=0.002434932427956787; "I2"=0.8105788195493853; "J2"=0.2109904491734014; "L2"=0.5103138085203796; "N2"=1.431523902812856; "O2"=3.278104424071699;
    "B3"=1.074822203913925; "D3"=0.1863123682828522; "E3"=0.1728898952792619; "F3"=1.313367131654047; "G3"=0.002438315577615514; "I3"=0.8233234718030724; "J3"=0.2078312979552095; "L3"=0.4702835370184459; "N3"=1.424807753880273; "O3"=3.244091665026531;
    "B4"=1.034072198655821; "D4"=0.1864608979753548; "E4"=0.172197003375393; "F4"=1.311466661671808; "G4"=0.002440506441020734; "I4"=0.8316012188112829; "J4"=0.2059032326143893; "L4"=0.4457670971699486; "N4"=1.421200269464322; "O4"=3.225330344706492;
    "B5"=1.017516894567962; "D5"=0.1865268445025379; "E5"=0.1719183987134691; "F5"=1.310893501922727; "G5"=0.002441427888173472; "I5"=0.8350881731822719; "J5"=0.2051205708363426; "L5"=0.4357927657334528; "N5"=1.419860627256938; "O5"=3.218218228765721;
    "B6"=1.014770999758042; "D6"=0.1865381227868887; "E6"=0.1718723649475162; "F6"=1.310810484271663; "G6"=0.002441582626811096; "I6"=0.835674044645768; "J6"=0.2049907972292502; "L6"=0.4341375402279084; "N6"=1.419646077118585; "O6"=3.217069461766471;
    "B7"=1.033848721075998; "D7"=0.186461765382905; "E7"=0.1721932307427814; "F7"=1.311458116923902; "G7"=0.002440518751888881; "I7"=0.8316477848082489; "J7"=0.2058926649106496; "L7"=0.4456325129207244; "N7"=1.421181673620353; "O7"=3.225232269665355;
    "B8"=1.118475448844549; "D8"=0.1861723416650634; "E8"=0.1736419588247209; "F8"=1.31608345802168; "G8"=0.002436075414241965; "I8"=0.8148792204003774; "J8"=0.2098988155353183; "L8"=0.4964988148595353; "N8"=1.42910138969367; "O8"=3.265935956914547;
    "B9"=1.285929011791382; "D9"=0.1857743057265608; "E9"=0.1765975237177813; "F9"=1.331515316039443; "G9"=0.002428259374766562; "I9"=0.785589177702156; "J9"=0.2178431354123092; "L9"=0.596718415990722; "N9"=1.448703326550188; "O9"=3.362629007510179;
    "B10"=1.409813326924734; "D10"=0.1855841105257809; "E10"=0.178836179211519; "F10"=1.346744252917347; "G10"=0.002423058378795317; "I10"=0.7662634171690179; "J10"=0.2237284554503773; "L10"=0.6706127810084297; "N10"=1.465556899515519; "O10"=3.444011823984511;
    "B11"=1.466345784256248; "D11"=0.1855195576819; "E11"=0.1798686441782493; "F11"=1.354520241602884; "G11"=0.002420808692244438; "I11"=0.7579486260310198; "J11"=0.2264153370622992; "L11"=0.7042817218397488; "N11"=1.47375044232075; "O11"=3.483293628065496;
    "B12"=1.487777335823182; "D12"=0.1854982529648481; "E12"=0.1802615868294879; "F12"=1.35758696293972; "G12"=0.002419973424324779; "I12"=0.7548686449254962; "J12"=0.2274340711109701; "L12"=0.7170384873236912; "N12"=1.476928324018786; "O12"=3.49849443538875;
    "B13"=1.483160626892641; "D13"=0.1855027019627649; "E13"=0.1801768728383379; "F13"=1.356921055889231; "G13"=0.002420152575457021; "I13"=0.7555289188900023; "J13"=0.2272146138065025; "L13"=0.714290786391814; "N13"=1.476240576567662; "O13"=3.495206176974136;
    "B14"=1.468108499720302; "D14"=0.1855177421068284; "E14"=0.1799009326303924; "F14"=1.354770094263273; "G14"=0.002420739641326528; "I14"=0.7576938576198504; "J14"=0.22649912411363; "L14"=0.7053310909584241; "N14"=1.474010385459835; "O14"=3.484537677847868;
    "B15"=1.458891710193257; "D15"=0.1855273629978313; "E15"=0.1797321662699503; "F15"=1.353468476448285; "G15"=0.002421101400811689; "I15"=0.7590288887942158; "J15"=0.2260610280977104; "L15"=0.6998439185535119; "N15"=1.472654100507171; "O15"=3.478045340045242;
    "B16"=1.406122228828508; "D16"=0.1855887696358529; "E16"=0.1787689839669326; "F16"=1.346253155348421; "G16"=0.002423207733590949; "I16"=0.766816408540155; "J16"=0.223553045325751; "L16"=0.6684134637324632; "N16"=1.465031982059187; "O16"=3.441490201483759;
    "B17"=1.373794093927017; "D17"=0.1856320536123697; "E17"=0.1781816719272236; "F17"=1.342044153483087; "G17"=0.002424529623697153; "I17"=0.771715949412604; "J17"=0.2220168625056189; "L17"=0.64914523560509; "N17"=1.46049055929069; "O17"=3.419644169009132;
    "B18"=1.355216541845721; "D18"=0.1856590182563771; "E18"=0.1778451950382447; "F18"=1.339703081461437; "G18"=0.002425300888680166; "I18"=0.7745788880657543; "J18"=0.221134202999032; "L18"=0.6380677870949967; "N18"=1.457928069296912; "O18"=3.407291627888128;
    "B19"=1.348929422007188; "D19"=0.1856685039262587; "E19"=0.1777314997736639; "F19"=1.338924141997609; "G19"=0.002425563908871449; "I19"=0.7755559298020778; "J19"=0.2208355097800876; "L19"=0.6343180552484569; "N19"=1.457068993776431; "O19"=3.403145785974743;
    "B20"=1.377233757800468; "D20"=0.1856272319794634; "E20"=0.1782440551048445; "F20"=1.342483945302632; "G20"=0.002424387773787359; "I20"=0.7711897420331137; "J20"=0.2221802981986016; "L20"=0.6511958458690117; "N20"=1.46096887065238; "O20"=3.421947696455618;
    "B21"=1.472529036593528; "D21"=0.1855132393822139; "E21"=0.1799819299186183; "F21"=1.355398568304707; "G21"=0.002420566754868619; "I21"=0.7570560985709882; "J21"=0.2267092472040062; "L21"=0.7079625844321811; "N21"=1.474663411308626; "O21"=3.487662432039599;
    "B22"=1.534948842073391; "D22"=0.1854570325560552; "E22"=0.1811291919901095; "F22"=1.364550914705305; "G22"=0.002418166448401963; "I22"=0.7482191435654073; "J22"=0.2296765292034308; "L22"=0.7451038056167363; "N22"=1.48405143836662; "O22"=3.532509173684673;
    "B23"=1.501622006886464; "D23"=0.1854853635028739; "E23"=0.1805158462773022; "F23"=1.359600949607838; "G23"=0.002419438692557518; "I23"=0.7528989361254492; "J23"=0.2280921995159915; "L23"=0.7252773275563129; "N23"=1.479000998859732; "O23"=3.5083997236747;
    "B24"=1.375678659716925; "D24"=0.1856294053597409; "E24"=0.1782158480019866; "F24"=1.342284870160512; "G24"=0.00242445186897287; "I24"=0.7714274967129722; "J24"=0.2221064073176606; "L24"=0.6502687645655101; "N24"=1.460752475169869; "O24"=3.420905626833303;
    "B25"=1.240473806810883; "D25"=0.1858639511096705; "E25"=0.1757859714884766; "F25"=1.326658199459771; "G25"=0.002430278335057542; "I25"=0.7931279190617433; "J25"=0.2178431354123092; "L25"=0.5695584373445968; "N25"=1.442967745587509; "O25"=3.334658785307766;
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}
